$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Update sheet "o_10": add a 5th column "evaluator_partial_correctness"
#    (copying the header style from D1) and update the existing
#    prompt/solution/llm_response/evaluator_response row for the
#    A-to-J question, plus its new partial-correctness score.
# ---------------------------------------------------------------------------
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws1.Range("E1").Value = "evaluator_partial_correctness"

$prompt_o10 = @'
 Given is the adjacency matrix for a unweighted undirected graph containing 10 nodes labelled A to J. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node O? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O
 A 0 1 0 1 0 0 0 0 0 0 0 0 1 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 1 0 0 0 1 0 0 1 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 L 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 M 1 0 0 0 0 0 0 0 0 0 0 0 0 1 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0

Solution: A -> M -> N -> O
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node J? Return the sequence of nodes in response.

   A B C D E F G H I J
 A 0 1 1 1 0 0 0 0 0 0
 B 1 0 0 0 0 0 0 0 0 0
 C 1 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 1 0 0 1 0 0
 E 0 0 0 1 0 1 1 0 0 0
 F 0 0 0 0 1 0 0 0 0 0
 G 0 0 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 0 0 1 0
 I 0 0 0 0 0 0 0 1 0 1
 J 0 0 0 0 0 0 0 0 1 0
    
'@

$ws1.Range("A2").Value = $prompt_o10
$ws1.Range("B2").Value = "A -> D -> H -> I -> J"
$ws1.Range("C2").Value = "The shortest path from node A to node J is: A -> D -> H -> I -> J"
$ws1.Range("D2").Value = "invalid input"
$ws1.Range("E2").Value = "5/5"

# ---------------------------------------------------------------------------
# 2) Add a new worksheet "o_20" right after "o_10" (duplicate "o_10" so it
#    inherits the same column layout/styles, then overwrite the data).
# ---------------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "o_20"

$prompt_o20 = @'
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 G 0 0 0 0 0 1 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 1 0 0 0 1 1 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 1 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 W 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0

Solution: A -> E -> F -> W -> Y
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node T? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 1 0 0 0 0 0 0 1 0 0 0 0 0
 G 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 1 0 0 1 1 1 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 1 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 1
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
    
'@

$ws2.Range("A2").Value = $prompt_o20
$ws2.Range("B2").Value = "A -> E -> F -> O -> P -> Q -> R -> T"
$ws2.Range("C2").Value = "The shortest path from node A to node T is: A -> E -> F -> O -> P -> Q -> R -> T"
$ws2.Range("D2").Value = "invalid input"
$ws2.Range("E2").Value = "4/4"

# ---------------------------------------------------------------------------
# 3) Add a new worksheet "o_20_jumbled" right after "o_20" (duplicate "o_20"
#    and overwrite the data with the jumbled-graph version).
# ---------------------------------------------------------------------------
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "o_20_jumbled"

$prompt_o20_jumbled = @'
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 G 0 0 0 0 0 1 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 1 0 0 0 1 1 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 1 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 W 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0

Solution: A -> E -> F -> W -> Y
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node T? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 1 0 0 1 1 1 0 1 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 1 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
'@

$ws3.Range("A2").Value = $prompt_o20_jumbled
$ws3.Range("B2").Value = "A -> B -> C -> E -> J -> K -> M -> N -> S -> T"
$ws3.Range("C2").Value = "The shortest path from node A to node T is: A -> B -> C -> E -> J -> K -> M -> N -> S -> T"
$ws3.Range("D2").Value = "invalid input"
$ws3.Range("E2").Value = "10/10"

# ---------------------------------------------------------------------------
# Leave "o_10" as the active sheet, matching the original workbook state.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
